# Fruta / hortaliza, semanal
# Insert two new data rows (new rows 112 and 113) into the daily price table,
# pushing the existing rows 112-217 down to 114-219.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at position 112, shifting everything below down.
$ws.Rows("112:113").Insert()

# Populate new row 112
$ws.Range("A112").Value = 6
$ws.Range("B112").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C112").Value = "Metropolitana"
$ws.Range("D112").Value = 44897
$ws.Range("E112").Value = 13
$ws.Range("F112").Value = "Fruta"
$ws.Range("G112").Value = 100101
$ws.Range("H112").Value = "Berries"
$ws.Range("I112").Value = 100101004
$ws.Range("J112").Value = "Frambuesa"
$ws.Range("K112").Value = "Sin especificar"
$ws.Range("L112").Value = "Especial"
$ws.Range("M112").Value = 200
$ws.Range("N112").Value = 9000
$ws.Range("O112").Value = 9000
$ws.Range("P112").Value = 9000
$ws.Range("Q112").Value = "$/bandeja 2 kilos"
$ws.Range("R112").Value = "Provincia de Curicó"
$ws.Range("S112").Value = 4500
$ws.Range("T112").Value = 2

# Populate new row 113
$ws.Range("A113").Value = 6
$ws.Range("B113").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C113").Value = "Metropolitana"
$ws.Range("D113").Value = 44897
$ws.Range("E113").Value = 13
$ws.Range("F113").Value = "Fruta"
$ws.Range("G113").Value = 100101
$ws.Range("H113").Value = "Berries"
$ws.Range("I113").Value = 100101004
$ws.Range("J113").Value = "Frambuesa"
$ws.Range("K113").Value = "Sin especificar"
$ws.Range("L113").Value = "Primera"
$ws.Range("M113").Value = 150
$ws.Range("N113").Value = 8000
$ws.Range("O113").Value = 8000
$ws.Range("P113").Value = 8000
$ws.Range("Q113").Value = "$/bandeja 2 kilos"
$ws.Range("R113").Value = "Provincia de Curicó"
$ws.Range("S113").Value = 4000
$ws.Range("T113").Value = 2
